$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2020" column (Q) by extending the existing year columns -----
# Column P (2019) is the template: copy its formatting across to Q for the
# whole data block (rows 3-34) and then fill in the 2020 figures.
$ws.Range("P3:P34").Copy()
$ws.Range("Q3:Q34").PasteSpecial(-4122)

# Header (row 4) - new year value
$ws.Range("Q4").Value = 2020

# Republic-wide total (row 5) and gender breakdown (rows 6-7)
$ws.Range("Q5").Value = 51
$ws.Range("Q6").Value = 29
$ws.Range("Q7").Value = 22

# Batken oblast (rows 8-10)
$ws.Range("Q8").Value = 5
$ws.Range("Q9").Value = 3
$ws.Range("Q10").Value = 2

# Djalal-Abad oblast (rows 11-13)
$ws.Range("Q11").Value = 15
$ws.Range("Q12").Value = 9
$ws.Range("Q13").Value = 5

# Ysyk-Kul oblast (rows 14-16) - no data for this year
$ws.Range("Q14").Value = "-"
$ws.Range("Q15").Value = "-"
$ws.Range("Q16").Value = "-"

# Naryn oblast (rows 17-19) - no data for this year
$ws.Range("Q17").Value = "-"
$ws.Range("Q18").Value = "-"
$ws.Range("Q19").Value = "-"

# Osh oblast (rows 20-22)
$ws.Range("Q20").Value = 7
$ws.Range("Q21").Value = 7
$ws.Range("Q22").Value = "-"

# Talas oblast (rows 23-25) - no data for this year
$ws.Range("Q23").Value = "-"
$ws.Range("Q24").Value = "-"
$ws.Range("Q25").Value = "-"

# Chui oblast (rows 26-28)
$ws.Range("Q26").Value = 24
$ws.Range("Q27").Value = 10
$ws.Range("Q28").Value = 14

# Bishkek city (rows 29-31) - no data for this year
$ws.Range("Q29").Value = "-"
$ws.Range("Q30").Value = "-"
$ws.Range("Q31").Value = "-"

# Osh city (rows 32-34) - no data for this year
$ws.Range("Q32").Value = "-"
$ws.Range("Q33").Value = "-"
$ws.Range("Q34").Value = "-"

# --- Restore the (arbitrary) active-cell selection recorded in the file -----
$ws.Range("K18").Select() | Out-Null
